# Weekly update to the Choclo (corn) price sheet:
# Insert a brand-new record as row 546, pushing the former rows 546..634
# down to 547..635 (dimension grows from A1:R634 to A1:R635).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 546, carrying formatting
# down from the row above (this is how Excel's own Rows.Insert behaves).
$ws.Rows.Item(546).Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(546, 1).Value  = 6
$ws.Cells.Item(546, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(546, 3).Value  = "Metropolitana"
$ws.Cells.Item(546, 4).Value  = 44474
$ws.Cells.Item(546, 5).Value  = 13
$ws.Cells.Item(546, 6).Value  = 100112024
$ws.Cells.Item(546, 7).Value  = "Choclo"
$ws.Cells.Item(546, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(546, 9).Value  = "Primera"
$ws.Cells.Item(546, 10).Value = 200
$ws.Cells.Item(546, 11).Value = 22000
$ws.Cells.Item(546, 12).Value = 23000
$ws.Cells.Item(546, 13).Value = 22600
$ws.Cells.Item(546, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(546, 15).Value = "Argentina"
$ws.Cells.Item(546, 16).Value = 452
$ws.Cells.Item(546, 17).Value = 50
$ws.Cells.Item(546, 18).Value = "Hortaliza"
